$wb = $excel.ActiveWorkbook

# --- Sheet ALC: 31 cell update(s) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3637.3333
$ws.Range("J76").Value = 3876.5
$ws.Range("L76").Value = 3876.5
$ws.Range("N76").Value = -4506.5
$ws.Range("H79").Value = 3637.3333
$ws.Range("J79").Value = 3876.5
$ws.Range("L79").Value = 3876.5
$ws.Range("N79").Value = -6060.5
$ws.Range("H98").Value = 4999.6665
$ws.Range("I98").Value = 4999.6665
$ws.Range("K98").Value = 4999.6665
$ws.Range("M98").Value = -3501.6665
$ws.Range("H113").Value = 4498
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 1254
$ws.Range("H122").Value = 4999.6665
$ws.Range("I122").Value = 4999.6665
$ws.Range("K122").Value = 14998.9995
$ws.Range("M122").Value = -12548.9995
$ws.Range("H137").Value = 11909.2
$ws.Range("I137").Value = 2549.1428
$ws.Range("K137").Value = 7647.428400000001
$ws.Range("M137").Value = -5097.428400000001
$ws.Range("H138").Value = 2337.677
$ws.Range("I138").Value = 1545.6364
$ws.Range("J138").Value = 2499.0186
$ws.Range("K138").Value = 4636.9092
$ws.Range("L138").Value = 7497.0558
$ws.Range("M138").Value = 503.0907999999999
$ws.Range("N138").Value = -17777.0558

# --- Sheet ARM: 23 cell update(s) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H32").Value = 228168.33
$ws.Range("I32").Value = 233060.55
$ws.Range("K32").Value = 233060.55
$ws.Range("M32").Value = -232773.55
$ws.Range("H61").Value = 15000
$ws.Range("I61").Value = 11666.667
$ws.Range("K61").Value = 11666.667
$ws.Range("M61").Value = -11454.667
$ws.Range("H132").Value = 534279.0600000001
$ws.Range("I132").Value = 627379.8
$ws.Range("J132").Value = 2274.7144
$ws.Range("K132").Value = 1882139.4
$ws.Range("L132").Value = 6824.1432
$ws.Range("M132").Value = -1879609.4
$ws.Range("N132").Value = -11884.1432
$ws.Range("H136").Value = 15000
$ws.Range("I136").Value = 11666.667
$ws.Range("K136").Value = 35000.001
$ws.Range("M136").Value = -32450.001

# --- Sheet CRP: 43 cell update(s) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2322.8286
$ws.Range("I31").Value = 2264.52
$ws.Range("J31").Value = 2468.6
$ws.Range("K31").Value = 2264.52
$ws.Range("L31").Value = 2468.6
$ws.Range("M31").Value = -1969.52
$ws.Range("N31").Value = -3058.6
$ws.Range("H34").Value = 2322.8286
$ws.Range("I34").Value = 2264.52
$ws.Range("J34").Value = 2468.6
$ws.Range("K34").Value = 2264.52
$ws.Range("L34").Value = 2468.6
$ws.Range("M34").Value = -2062.52
$ws.Range("N34").Value = -2872.6
$ws.Range("H105").Value = 12646.223
$ws.Range("I105").Value = 14988
$ws.Range("K105").Value = 14988
$ws.Range("M105").Value = -13241
$ws.Range("H106").Value = 40000
$ws.Range("J106").Value = 40000
$ws.Range("L106").Value = 40000
$ws.Range("N106").Value = -42524
$ws.Range("H107").Value = 542.2083
$ws.Range("I107").Value = 514.9474
$ws.Range("J107").Value = 645.8
$ws.Range("K107").Value = 514.9474
$ws.Range("L107").Value = 645.8
$ws.Range("M107").Value = 1405.0526
$ws.Range("N107").Value = -4485.8
$ws.Range("H132").Value = 2310.05
$ws.Range("I132").Value = 1659.2
$ws.Range("J132").Value = 4262.6
$ws.Range("K132").Value = 4977.6
$ws.Range("L132").Value = 12787.8
$ws.Range("M132").Value = -2447.6
$ws.Range("N132").Value = -17847.8
$ws.Range("H141").Value = 353263.62
$ws.Range("I141").Value = 25000
$ws.Range("J141").Value = 400158.44
$ws.Range("K141").Value = 25000
$ws.Range("L141").Value = 400158.44
$ws.Range("M141").Value = -19820
$ws.Range("N141").Value = -410518.44

# --- Sheet CUL: 15 cell update(s) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 11991
$ws.Range("I56").Value = 11991
$ws.Range("K56").Value = 11991
$ws.Range("M56").Value = -11461
$ws.Range("H108").Value = 506.875
$ws.Range("I108").Value = 506.875
$ws.Range("K108").Value = 1520.625
$ws.Range("M108").Value = 1359.375
$ws.Range("H132").Value = 768.53845
$ws.Range("I132").Value = 563.4286
$ws.Range("J132").Value = 1007.8333
$ws.Range("K132").Value = 5070.8574
$ws.Range("L132").Value = 9070.4997
$ws.Range("M132").Value = -2540.8574
$ws.Range("N132").Value = -14130.4997

# --- Sheet GSM: 28 cell update(s) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1972.5
$ws.Range("I80").Value = 1945
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 1945
$ws.Range("L80").Value = 2000
$ws.Range("M80").Value = -947
$ws.Range("N80").Value = -3996
$ws.Range("H83").Value = 1972.5
$ws.Range("I83").Value = 1945
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 9725
$ws.Range("L83").Value = 10000
$ws.Range("M83").Value = -4733
$ws.Range("N83").Value = -19984
$ws.Range("H102").Value = 1672.5161
$ws.Range("I102").Value = 1575.8928
$ws.Range("J102").Value = 2574.3333
$ws.Range("K102").Value = 1575.8928
$ws.Range("L102").Value = 2574.3333
$ws.Range("M102").Value = 46.10719999999992
$ws.Range("N102").Value = -5818.3333
$ws.Range("H132").Value = 10370.296
$ws.Range("I132").Value = 12892.605
$ws.Range("J132").Value = 4379.8125
$ws.Range("K132").Value = 38677.815
$ws.Range("L132").Value = 13139.4375
$ws.Range("M132").Value = -36147.815
$ws.Range("N132").Value = -18199.4375

# --- Sheet LTW: 37 cell update(s) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5535.6665
$ws.Range("I7").Value = 3966.4
$ws.Range("K7").Value = 3966.4
$ws.Range("M7").Value = -3854.4
$ws.Range("H40").Value = 3159.9
$ws.Range("I40").Value = 1934
$ws.Range("J40").Value = 4998.75
$ws.Range("K40").Value = 1934
$ws.Range("L40").Value = 4998.75
$ws.Range("M40").Value = -1798
$ws.Range("N40").Value = -5270.75
$ws.Range("H82").Value = 2851.5715
$ws.Range("I82").Value = 3592.8
$ws.Range("K82").Value = 3592.8
$ws.Range("M82").Value = -3231.8
$ws.Range("H85").Value = 2851.5715
$ws.Range("I85").Value = 3592.8
$ws.Range("K85").Value = 3592.8
$ws.Range("M85").Value = -2344.8
$ws.Range("H126").Value = 5535.6665
$ws.Range("I126").Value = 3966.4
$ws.Range("K126").Value = 11899.2
$ws.Range("M126").Value = -9429.200000000001
$ws.Range("H132").Value = 1430461.9
$ws.Range("I132").Value = 1726044.6
$ws.Range("J132").Value = 1812
$ws.Range("K132").Value = 5178133.800000001
$ws.Range("L132").Value = 5436
$ws.Range("M132").Value = -5175603.800000001
$ws.Range("N132").Value = -10496
$ws.Range("H136").Value = 8930.049999999999
$ws.Range("I136").Value = 4404.636
$ws.Range("J136").Value = 14461.111
$ws.Range("K136").Value = 13213.908
$ws.Range("L136").Value = 43383.333
$ws.Range("M136").Value = -10663.908
$ws.Range("N136").Value = -48483.333

# --- Sheet WVR: 25 cell update(s) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2816.611
$ws.Range("I81").Value = 2484.5386
$ws.Range("J81").Value = 3680
$ws.Range("K81").Value = 4969.0772
$ws.Range("L81").Value = 7360
$ws.Range("M81").Value = -3908.0772
$ws.Range("N81").Value = -9482
$ws.Range("H84").Value = 2816.611
$ws.Range("I84").Value = 2484.5386
$ws.Range("J84").Value = 3680
$ws.Range("K84").Value = 24845.386
$ws.Range("L84").Value = 36800
$ws.Range("M84").Value = -19541.386
$ws.Range("N84").Value = -47408
$ws.Range("H122").Value = 2112.0293
$ws.Range("I122").Value = 1955.409
$ws.Range("J122").Value = 2399.1667
$ws.Range("K122").Value = 5866.227000000001
$ws.Range("L122").Value = 7197.500100000001
$ws.Range("M122").Value = -3416.227000000001
$ws.Range("N122").Value = -12097.5001
$ws.Range("H136").Value = 1559.88
$ws.Range("I136").Value = 1526.2106
$ws.Range("K136").Value = 4578.6318
$ws.Range("M136").Value = -2028.6318
